$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Insert four new rows above the current row 33 (Biopsies_yesno). This
#    shifts the existing rows 33-42 down to 37-46, exactly like the target
#    diff (old row 33 -> new row 37, ... old row 42 -> new row 46).
# ---------------------------------------------------------------------------
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(33).Insert()

# ---------------------------------------------------------------------------
# 2. Populate the four new rows (BBD_year1 .. BBD_year4 valid-range rules).
# ---------------------------------------------------------------------------
$ws.Cells.Item(33, 1).Value = "BBD_year1"
$ws.Cells.Item(33, 2).Value = 777
$ws.Cells.Item(33, 3).Value = "valid.changes"
$ws.Cells.Item(33, 12).Value = 7777
$ws.Cells.Item(33, 13).Value = " BBD_year1 changed to 7777 from 777 to match data dictionary"

$ws.Cells.Item(34, 1).Value = "BBD_year2"
$ws.Cells.Item(34, 2).Value = 777
$ws.Cells.Item(34, 3).Value = "valid.changes"
$ws.Cells.Item(34, 12).Value = 7777
$ws.Cells.Item(34, 13).Value = " BBD_year2 changed to 7777 from 777 to match data dictionary"

$ws.Cells.Item(35, 1).Value = "BBD_year3"
$ws.Cells.Item(35, 2).Value = 777
$ws.Cells.Item(35, 3).Value = "valid.changes"
$ws.Cells.Item(35, 12).Value = 7777
$ws.Cells.Item(35, 13).Value = " BBD_year3 changed to 7777 from 777 to match data dictionary"

$ws.Cells.Item(36, 1).Value = "BBD_year4"
$ws.Cells.Item(36, 2).Value = 777
$ws.Cells.Item(36, 3).Value = "valid.changes"
$ws.Cells.Item(36, 12).Value = 7777
$ws.Cells.Item(36, 13).Value = " BBD_year4 changed to 7777 from 777 to match data dictionary"

# Column D on these four rows stays blank (matches the source diff, which
# only carries the style forward) - nothing to write.

# ---------------------------------------------------------------------------
# 3. Append two new rows (47-48) describing the "lastfup" correction rule.
# ---------------------------------------------------------------------------
$ws.Cells.Item(47, 1).Value = "lastfup"
$ws.Cells.Item(47, 2).Value = 777
$ws.Cells.Item(47, 12).Value = 7777
$ws.Cells.Item(47, 13).Value = "lastfup changed to 7777 from 777 to match data dictionary"

$ws.Cells.Item(48, 1).Value = "lastfup"
$ws.Cells.Item(48, 2).Value = 888
$ws.Cells.Item(48, 12).Value = 8888
$ws.Cells.Item(48, 13).Value = "lastfup changed to 8888 from 888 to match data dictionary"

# ---------------------------------------------------------------------------
# 4. Restore the view so the new bottom of the sheet is visible, matching
#    the saved workbook view in the target file.
# ---------------------------------------------------------------------------
$ws.Range("A37").Select()
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("M49").Select()
